$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$s.Shapes.Item("Table 3").Delete()
